# Applies updated loading-percent results for the 380 kV case (Case_2_247)
# Source sheet: Code/Results/Cases/Case_2_247/res_line/loading_percent.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 11.67584708394852
$rowData[0,1] = 8.691826167969872
$rowData[0,2] = 0
$rowData[0,3] = 25.69665789997646
$rowData[0,4] = 37.87022561806203
$rowData[0,5] = 18.43273599044
$rowData[0,6] = 11.59114964640461
$rowData[0,7] = 15.76462976287502
$rowData[0,8] = 7.20562058375808
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.62113110960501
$rowData[0,13] = 16.26201429166088
$ws.Range("B2:O2").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 11.05862818438093
$rowData[0,1] = 8.274085650095671
$rowData[0,2] = 0
$rowData[0,3] = 25.43582716291737
$rowData[0,4] = 37.62526164352224
$rowData[0,5] = 18.39976966432056
$rowData[0,6] = 11.63436161410086
$rowData[0,7] = 15.87041583954898
$rowData[0,8] = 7.226373427101731
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.63009504965982
$rowData[0,13] = 16.31917394490299
$ws.Range("B3:O3").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 10.66137984034244
$rowData[0,1] = 8.005132440391442
$rowData[0,2] = 0
$rowData[0,3] = 25.28038776866602
$rowData[0,4] = 37.48526532848666
$rowData[0,5] = 18.38870407784368
$rowData[0,6] = 11.66310710580771
$rowData[0,7] = 15.93912663380821
$rowData[0,8] = 7.239917429199235
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.63725270982864
$rowData[0,13] = 16.35888258343732
$ws.Range("B4:O4").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 10.49505142632584
$rowData[0,1] = 7.892486368674409
$rowData[0,2] = 0
$rowData[0,3] = 25.21829406996816
$rowData[0,4] = 37.43088411252668
$rowData[0,5] = 18.38649968433593
$rowData[0,6] = 11.67537696093435
$rowData[0,7] = 15.96807179103775
$rowData[0,8] = 7.24563865659809
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.6405866090926
$rowData[0,13] = 16.37621914577282
$ws.Range("B5:O5").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 10.46716906556263
$rowData[0,1] = 7.873600494666776
$rowData[0,2] = 0
$rowData[0,3] = 25.20806073001049
$rowData[0,4] = 37.42201668549973
$rowData[0,5] = 18.38627271497262
$rowData[0,6] = 11.67744791199932
$rowData[0,7] = 15.97293518761436
$rowData[0,8] = 7.246600868867765
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.64116542738424
$rowData[0,13] = 16.37916748046648
$ws.Range("B6:O6").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 10.65915446381834
$rowData[0,1] = 8.003625460244916
$rowData[0,2] = 0
$rowData[0,3] = 25.2795452107243
$rowData[0,4] = 37.48452105764282
$rowData[0,5] = 18.38866502205487
$rowData[0,6] = 11.66327033173138
$rowData[0,7] = 15.93951317215814
$rowData[0,8] = 7.239993769570926
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.6372959815585
$rowData[0,13] = 16.35911172052739
$ws.Range("B7:O7").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 11.46690333427005
$rowData[0,1] = 8.550422930155216
$rowData[0,2] = 0
$rowData[0,3] = 25.60578489406827
$rowData[0,4] = 37.78363246193616
$rowData[0,5] = 18.41946442236219
$rowData[0,6] = 11.60558937918054
$rowData[0,7] = 15.80032490355662
$rowData[0,8] = 7.212609987149356
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.62387939422112
$rowData[0,13] = 16.2807627823564
$ws.Range("B8:O8").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 12.90068338000443
$rowData[0,1] = 9.520711403837282
$rowData[0,2] = 0
$rowData[0,3] = 26.27998597412625
$rowData[0,4] = 38.4503764128124
$rowData[0,5] = 18.55259621490075
$rowData[0,6] = 11.51006824991942
$rowData[0,7] = 15.55719111875405
$rowData[0,8] = 7.165255581217536
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.61063625083421
$rowData[0,13] = 16.16392638422517
$ws.Range("B9:O9").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.85702428042298
$rowData[0,1] = 10.16807736594849
$rowData[0,2] = 0
$rowData[0,3] = 26.79225786204916
$rowData[0,4] = 38.98555688868989
$rowData[0,5] = 18.69439191939657
$rowData[0,6] = 11.45065447836888
$rowData[0,7] = 15.3967341497446
$rowData[0,8] = 7.134310545894484
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60879738222842
$rowData[0,13] = 16.10079642436631
$ws.Range("B10:O10").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.27020546489841
$rowData[0,1] = 10.44785306704437
$rowData[0,2] = 0
$rowData[0,3] = 27.0280937984443
$rowData[0,4] = 39.23802794409785
$rowData[0,5] = 18.76829324772744
$rowData[0,6] = 11.42597222154727
$rowData[0,7] = 15.32768551686563
$rowData[0,8] = 7.121063450595573
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60965730724246
$rowData[0,13] = 16.07706149905112
$ws.Range("B11:O11").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.42347265252013
$rowData[0,1] = 10.55164910684147
$rowData[0,2] = 0
$rowData[0,3] = 27.11772649748585
$rowData[0,4] = 39.33485339180417
$rowData[0,5] = 18.79761043949016
$rowData[0,6] = 11.41696372493086
$rowData[0,7] = 15.30210617324477
$rowData[0,8] = 7.116166133513185
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.61022533133768
$rowData[0,13] = 16.068794231533
$ws.Range("B12:O12").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.39060665195794
$rowData[0,1] = 10.52939080430773
$rowData[0,2] = 0
$rowData[0,3] = 27.09840916684849
$rowData[0,4] = 39.31394727204587
$rowData[0,5] = 18.7912375408475
$rowData[0,6] = 11.41888881314211
$rowData[0,7] = 15.30758987358609
$rowData[0,8] = 7.117215567018715
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.61009224227809
$rowData[0,13] = 16.07054262374609
$ws.Range("B13:O13").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.28287914958275
$rowData[0,1] = 10.45643566377451
$rowData[0,2] = 0
$rowData[0,3] = 27.03546182037925
$rowData[0,4] = 39.24596982361849
$rowData[0,5] = 18.77067860970602
$rowData[0,6] = 11.42522430530367
$rowData[0,7] = 15.32556970331
$rowData[0,8] = 7.120658160584424
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.6096991902134
$rowData[0,13] = 16.07636688124035
$ws.Range("B14:O14").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.21647544177542
$rowData[0,1] = 10.41146780923965
$rowData[0,2] = 0
$rowData[0,3] = 26.99694505401671
$rowData[0,4] = 39.20448826859123
$rowData[0,5] = 18.75825855878567
$rowData[0,6] = 11.4291490387442
$rowData[0,7] = 15.33665684545067
$rowData[0,8] = 7.122782346931102
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60948995349052
$rowData[0,13] = 16.08002837015874
$ws.Range("B15:O15").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.82957730318165
$rowData[0,1] = 10.14949429788979
$rowData[0,2] = 0
$rowData[0,3] = 26.77689528711235
$rowData[0,4] = 38.9692322583772
$rowData[0,5] = 18.68974985575192
$rowData[0,6] = 11.45231479178453
$rowData[0,7] = 15.40132606448096
$rowData[0,8] = 7.135192951813504
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60877518906271
$rowData[0,13] = 16.10244818573646
$ws.Range("B16:O16").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.58658776533429
$rowData[0,1] = 9.984987662103832
$rowData[0,2] = 0
$rowData[0,3] = 26.64256512176369
$rowData[0,4] = 38.82716584701901
$rowData[0,5] = 18.65011657680043
$rowData[0,6] = 11.46712756129744
$rowData[0,7] = 15.44200918724711
$rowData[0,8] = 7.143018829159531
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60877011381109
$rowData[0,13] = 16.11748131320738
$ws.Range("B17:O17").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.44477105366318
$rowData[0,1] = 9.888984322945548
$rowData[0,2] = 0
$rowData[0,3] = 26.56557014727318
$rowData[0,4] = 38.7463062563566
$rowData[0,5] = 18.62820649438938
$rowData[0,6] = 11.47586817920996
$rowData[0,7] = 15.46578022646154
$rowData[0,8] = 7.147598202392663
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60892696803933
$rowData[0,13] = 16.1265967894006
$ws.Range("B18:O18").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.39640306832276
$rowData[0,1] = 9.856242746874489
$rowData[0,2] = 0
$rowData[0,3] = 26.53954935186754
$rowData[0,4] = 38.7190774138485
$rowData[0,5] = 18.62094079766269
$rowData[0,6] = 11.47886548169066
$rowData[0,7] = 15.47389242239925
$rowData[0,8] = 7.149162127864031
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60900756386337
$rowData[0,13] = 16.12976352482032
$ws.Range("B19:O19").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.61266756834409
$rowData[0,1] = 10.00264311907654
$rowData[0,2] = 0
$rowData[0,3] = 26.65683759251326
$rowData[0,4] = 38.84220129438008
$rowData[0,5] = 18.65424404291228
$rowData[0,6] = 11.46552786820565
$rowData[0,7] = 15.43763997835409
$rowData[0,8] = 7.14217766643896
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60875412921331
$rowData[0,13] = 16.11583245980854
$ws.Range("B20:O20").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.31460841257469
$rowData[0,1] = 10.47792290964393
$rowData[0,2] = 0
$rowData[0,3] = 27.05394271839234
$rowData[0,4] = 39.26590395306286
$rowData[0,5] = 18.77668127528668
$rowData[0,6] = 11.423354233733
$rowData[0,7] = 15.32027317325418
$rowData[0,8] = 7.119643758040313
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.6098080734591
$rowData[0,13] = 16.07463656700416
$ws.Range("B21:O21").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.75472453858812
$rowData[0,1] = 10.77600911892787
$rowData[0,2] = 0
$rowData[0,3] = 27.31534345469975
$rowData[0,4] = 39.54989307646259
$rowData[0,5] = 18.8644558386412
$rowData[0,6] = 11.39776265824441
$rowData[0,7] = 15.24687759750345
$rowData[0,8] = 7.10561045580351
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.61190895219396
$rowData[0,13] = 16.05191474180564
$ws.Range("B22:O22").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 14.52154642235479
$rowData[0,1] = 10.6180712664049
$rowData[0,2] = 0
$rowData[0,3] = 27.1756829113495
$rowData[0,4] = 39.39770115462044
$rowData[0,5] = 18.81690652881561
$rowData[0,6] = 11.41124069805449
$rowData[0,7] = 15.28574702948271
$rowData[0,8] = 7.113036890266447
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.61065899921874
$rowData[0,13] = 16.06365605566998
$ws.Range("B23:O23").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 13.60088348492037
$rowData[0,1] = 9.994665527006443
$rowData[0,2] = 0
$rowData[0,3] = 26.65038427616429
$rowData[0,4] = 38.83540122011128
$rowData[0,5] = 18.6523752865733
$rowData[0,6] = 11.46625039001232
$rowData[0,7] = 15.43961410853449
$rowData[0,8] = 7.142557706355023
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.60876085804801
$rowData[0,13] = 16.11657643418781
$ws.Range("B24:O24").Value = $rowData

$rowData = New-Object 'object[,]' 1,14
$rowData[0,0] = 12.52952813246416
$rowData[0,1] = 9.269525948514623
$rowData[0,2] = 0
$rowData[0,3] = 26.09432314391372
$rowData[0,4] = 38.26177153149801
$rowData[0,5] = 18.50881565714725
$rowData[0,6] = 11.53402170135608
$rowData[0,7] = 15.61977398548952
$rowData[0,8] = 7.177389201603762
$rowData[0,9] = 0
$rowData[0,10] = 0
$rowData[0,11] = 0
$rowData[0,12] = 15.61282750599472
$rowData[0,13] = 16.19156440078801
$ws.Range("B25:O25").Value = $rowData
